# 2-1-1.xlsx "Add files via upload" change
#
# Updates the "2. Data reporter" contact block (Organization, Contact
# person, email, phone, website) to the new National Statistical
# Committee contact and moves the active selection to B8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the contact-detail cells. The order below matters: it reproduces
# the insertion order of the underlying shared strings (website, contact
# person, email, phone, organization) produced by the original edit.
$ws.Range("B10").Value = "www.stat.gov.kg"
$ws.Range("B7").Value = "Kalymbetova Yryskan"
$ws.Range("B8").Value = "yryskan.kalymbetova@gmail.com"
$ws.Range("B9").Value = "(0312) 32 46 55"
$ws.Range("B6").Value = " National Statistical Committee of the Kyrgyz Republic (Department of Household Statistics)"

# The organization (B6) and phone (B9) cells were retyped directly in the
# source app, which re-stamped their font explicitly (still Calibri, same
# look) -- reproduce that so those two cells pick up fresh style entries,
# same as the authoring app did.
$ws.Range("B6").Font.Name = "Calibri"
$ws.Range("B9").Font.Name = "Calibri"

# Move the active selection to B8, matching the saved view state.
$ws.Range("B8").Select()
